$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed data values (rows 2-34) to reflect corrected weather cluster statistics
$ws.Range("E2").Value = 0
$ws.Range("G2").Value = 0.03548728813559317
$ws.Range("J2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("G3").Value = 0.06091101694915243
$ws.Range("H3").Value = 0.09234234234234238
$ws.Range("J3").Value = 0
$ws.Range("B4").Value = 0.112189859762675
$ws.Range("C4").Value = 0.07479045776918121
$ws.Range("D4").Value = 0.01621895590471363
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("G5").Value = 0.178495762711864
$ws.Range("H5").Value = 0.2972972972972971
$ws.Range("I5").Value = 0.1311672683513836
$ws.Range("J5").Value = 0.2965144820814914
$ws.Range("E6").Value = 0
$ws.Range("G6").Value = 0.0233050847457627
$ws.Range("H6").Value = 0.002252252252252252
$ws.Range("J6").Value = 0
$ws.Range("C7").Value = 0.02965828497743389
$ws.Range("B8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.4265672950460376
$ws.Range("F8").Value = 0.1254137550570053
$ws.Range("G8").Value = 0.07944915254237273
$ws.Range("H8").Value = 0.2252252252252255
$ws.Range("J8").Value = 0.3210603829160519
$ws.Range("C9").Value = 0.02578981302385556
$ws.Range("B10").Value = 0.04962243797195242
$ws.Range("C10").Value = 0.0006447453255963894
$ws.Range("D10").Value = 0.1449569183983778
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("K10").Value = 0.2480857580398166
$ws.Range("B12").Value = 0.007551240560949298
$ws.Range("C12").Value = 0.03546099290780141
$ws.Range("G12").Value = 0
$ws.Range("I12").Value = 0.002406738868832732
$ws.Range("E13").Value = 0
$ws.Range("G13").Value = 0.04555084745762705
$ws.Range("H13").Value = 0.004504504504504504
$ws.Range("J13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("G14").Value = 0.185381355932203
$ws.Range("I14").Value = 0.07581227436823101
$ws.Range("J14").Value = 0.1143838978890522
$ws.Range("D15").Value = 0.01571211353269133
$ws.Range("F15").Value = 0
$ws.Range("K15").Value = 0.03981623277182237
$ws.Range("B16").Value = 0.03451995685005387
$ws.Range("C16").Value = 0.07994842037395229
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 0.006355932203389832
$ws.Range("I16").Value = 0.2767749699157642
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = 0.004766949152542373
$ws.Range("J17").Value = 0
$ws.Range("B18").Value = 0.002157497303128371
$ws.Range("C18").Value = 0.001934235976789168
$ws.Range("G18").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 0.09004237288135575
$ws.Range("H19").Value = 0.002252252252252252
$ws.Range("J19").Value = 0.009327442317133039
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 0.02171610169491524
$ws.Range("J20").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("G21").Value = 0.04184322033898299
$ws.Range("H21").Value = 0.006756756756756757
$ws.Range("J21").Value = 0
$ws.Range("D22").Value = 0.05879371515458702
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 0.1244703389830506
$ws.Range("H23").Value = 0.123873873873874
$ws.Range("J23").Value = 0
$ws.Range("B24").Value = 0.006472491909385113
$ws.Range("C24").Value = 0.0006447453255963894
$ws.Range("G24").Value = 0
$ws.Range("B32").Value = 0.1359223300970873
$ws.Range("C32").Value = 0.1328175370728559
$ws.Range("D32").Value = 0.03395843892549421
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("I32").Value = 0.1817087845968713
$ws.Range("B33").Value = 0.2211434735706598
$ws.Range("C33").Value = 0.05609284332688597
$ws.Range("D33").Value = 0.00456158134820071
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("I33").Value = 0.004813477737665463
$ws.Range("K33").Value = 0.1592649310872894
$ws.Range("E34").Value = 0
$ws.Range("G34").Value = 0.003177966101694915

# Remove obsolete "Joint regime area" rows 36-40 (no longer part of the output)
$ws.Range("A36:K40").EntireRow.Delete()

